$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7644227147102356
$ws.Range("B1").Value = 3.929542779922485
$ws.Range("C1").Value = 2.576885938644409
$ws.Range("D1").Value = 2.151169776916504
$ws.Range("E1").Value = 1.963559627532959
